# Auto-generated edit script for "liste des requetes.xlsx"
# Applies the #131 query-rename + new-rows update to Feuil1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Update changed cells in existing rows (code-name column, plus a couple of minor text fixes) ---
$ws.Range("C2").Value = 'REALMS_UPD_IPBANNED_DEBANAUTOIP'
$ws.Range("C3").Value = 'REALMS_SEL_IPBANNED_INFOSSURIPBANNIES'
$ws.Range("C4").Value = 'REALMS_INS_ACCOUNTBANNED_AUTOBANCOMPTEAUTH'
$ws.Range("C5").Value = 'REALMS_UPD_ACCOUNT_MAJERREURSAUTH'
$ws.Range("C6").Value = 'REALMS_SEL_ACCOUNT_RECUPINFOSCOMPTE'
$ws.Range("E6").Value = 'SELECT sha_pass_hash,util_numero,locked,last_ip,failed_logins FROM account WHERE username = $1'
$ws.Range("C7").Value = 'REALMS_INS_AVERTISSEMENTS_AJOUTAVERTISSEMENTS'
$ws.Range("C8").Value = 'REALMS_UPD_ACCOUNT_MAJAVERTISSEMENTSCOMPTE'
$ws.Range("C9").Value = 'REALMS_INS_ACCOUNTBANNED_AUTOBANCOMPTEPOURAVERTISSEMENTS'
$ws.Range("D9").Value = 'AutoBan pour nombre d''avertissements'
$ws.Range("C10").Value = 'REALMS_UPD_ACCOUNTBANNED_DEBANCOMPTE'
$ws.Range("C11").Value = 'REALMS_SEL_LISTESERVEUR_RECUPLISTESERVEUR'
$ws.Range("C12").Value = 'REALMS_UPD_ACCOUNTACCESS_MAJGMLEVELCOMPTE'
$ws.Range("C13").Value = 'REALMS_UPD_IPBANNED_DEBANIP'
$ws.Range("C14").Value = 'REALMS_UPD_ACCOUNTBANNED_DEBANAUTOCOMPTE'

# --- Append new rows 16-24 ---
$ws.Range("A16").Value = 'liste_serveur'
$ws.Range("B16").Value = 'Update'
$ws.Range("C16").Value = 'REALMS_UPD_LISTESERVEUR_MAJCOMPTEONLINE'
$ws.Range("D16").Value = 'maj du nombre compte online dans la liste des serveurs'
$ws.Range("E16").Value = 'UPDATE Liste_serveur SET serveur_n_online = $1 WHERE id_serveur = $2'

$ws.Range("A17").Value = 'liste_serveur'
$ws.Range("B17").Value = 'Insert'
$ws.Range("C17").Value = 'REALMS_INS_LISTESERVEUR_CREATIONSERVEUR'
$ws.Range("D17").Value = 'creation d''un nouveau serveur'
$ws.Range("E17").Value = 'INSERT INTO Liste_serveur VALUES ($1, $2, $3, $4, ''0'', $5)'

$ws.Range("A18").Value = 'account_banned'
$ws.Range("B18").Value = 'Insert'
$ws.Range("C18").Value = 'REALMS_INS_ACCOUNTBANNED_BANCOMPTE'
$ws.Range("D18").Value = 'Ban d''un compte'
$ws.Range("E18").Value = 'INSERT INTO account_banned VALUES ($1, $2, $3, true, $4, $5'

$ws.Range("A19").Value = 'ip_banned'
$ws.Range("B19").Value = 'Insert'
$ws.Range("C19").Value = 'REALMS_INS_IPBANNED_BANIP'
$ws.Range("D19").Value = 'Ban d''une ip'
$ws.Range("E19").Value = 'INSERT INTO ip_banned VALUES ($1, $2, $3, $4, $5, true)'

$ws.Range("A20").Value = 'account'
$ws.Range("B20").Value = 'Update'
$ws.Range("C20").Value = 'REALMS_UPD_ACCOUNT_LOGCONNEXIONCOMPTE'
$ws.Range("D20").Value = 'Activation connexion compte'
$ws.Range("E20").Value = 'UPDATE account SET last_ip = $1, failed_logins = ''0'', last_login = now(), online = true WHERE Util_numero = $2'

$ws.Range("A21").Value = 'account'
$ws.Range("B21").Value = 'Update'
$ws.Range("C21").Value = 'REALMS_UPD_ACCOUNT_MAJMAIL'
$ws.Range("D21").Value = 'Maj mail associé au compte'
$ws.Range("E21").Value = 'UPDATE account SET email = $1 WHERE Util_numero = $2'

$ws.Range("A22").Value = 'account'
$ws.Range("B22").Value = 'Update'
$ws.Range("C22").Value = 'REALMS_UPD_ACCOUNT_MAJPASS'
$ws.Range("D22").Value = 'Maj pass compte'
$ws.Range("E22").Value = 'UPDATE account SET sha_pass_hash = $1 WHERE Util_numero = $2'

$ws.Range("A23").Value = 'account_access'
$ws.Range("B23").Value = 'Insert'
$ws.Range("C23").Value = 'REALMS_INS_ACCOUNTACCESS_CREATIONACCESS'
$ws.Range("D23").Value = 'Creation de l''access d''un compte'
$ws.Range("E23").Value = 'INSERT INTO account_access VALUES ($1, ''0'', SELECT MAX(Util_numero) from account, $4, $5, true)'

$ws.Range("A24").Value = 'Avertissements'
$ws.Range("B24").Value = 'Select'
$ws.Range("C24").Value = 'REALMS_SEL_AVERTISSEMENTS_LISTEAVERTOS'
$ws.Range("D24").Value = 'Liste des avertissements d''un compte'
$ws.Range("E24").Value = 'SELECT avertissement_date, avertissements_raison, username from Avertissements, account WHERE avertissements_numero_util = $1 and avertissements.avertissements_id_gm = account.util_numero)'

# --- Refresh the AutoFilter to cover the new extent ---
$ws.AutoFilterMode = $false
$ws.Range("A1:E24").AutoFilter()

# --- Update the hidden _FilterDatabase defined name to match ---
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Feuil1!`$A`$1:`$E`$24"

# --- Column width follow-up (C widened to fit the longer REALMS_* names, D resized) ---
$ws.Columns.Item(3).ColumnWidth = 67.7
$ws.Columns.Item(4).ColumnWidth = 43.7

# --- Selection ends on C25 (next empty "Nom requete" cell), view scrolled back to A1 ---
$ws.Range("A1").Select()
$ws.Range("C25").Select()
